$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D cells contain numeric-looking text that must remain as text (inline string)
# in the original workbook. We temporarily force a text number format so Excel
# does not silently convert the assigned string into a numeric value, then restore
# the cell style back to the workbook default ("Normal") to avoid leftover formatting.
$priceUpdates = @{
    "D2" = "247.96"
    "D3" = "21.70"
    "D4" = "5.434"
    "D5" = "0.05690"
    "D6" = "3.380"
    "D7" = "0.8054"
    "D8" = "1.034"
    "D9" = "0.1466"
    "D10" = "0.07780"
    "D12" = "0.03046"
    "D13" = "0.09267"
    "D14" = "3.585"
    "D15" = "0.001642"
    "D16" = "0.04702"
    "D17" = "0.0005861"
    "D18" = "0.006350"
    "D19" = "0.005050"
    "D20" = "0.001045"
    "D21" = "0.0001500"
    "D22" = "0.0003200"
    "D23" = "3.771"
    "D24" = "6.427"
    "D25" = "2.143"
    "D26" = "0.3263"
    "D27" = "0.1301"
    "D40" = "0.04114"
    "D41" = "0.006981"
    "D42" = "0.1046"
    "D43" = "0.003196"
    "D44" = "0.008100"
    "D45" = "0.00005892"
    "D46" = "0.00000000750"
    "D47" = "0.0005501"
    "D48" = "0.6825"
    "D49" = "0.009172"
    "D50" = "0.00002100"
    "D51" = "0.01010"
}

foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
    $cell.Style = "Normal"
}

# Column E text updates (plain text, no numeric-format workaround needed)
$ws.Range("E19").Value = "18HotbitTokenHTB"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
